$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")
$helper.Interior.ThemeColor = 5
$helper.Interior.TintAndShade = 0.6
$helper.Copy()
$target = $ws.Range("A34:C35")
$target.PasteSpecial(-4122)
$helper.Clear()
Write-Host "done"
